$wb = $excel.ActiveWorkbook

# ---------- Sheet "Entities" ----------
$ws1 = $wb.Worksheets.Item("Entities")

# Data edits
$ws1.Range("J2").Value = 6
$ws1.Range("G3").Value = "Scavenge for old keys."
$ws1.Range("I3").Value = 10201
$ws1.Range("K3").Value = "Use old key"
$ws1.Range("M3").Value = 10201

# Column width adjustments (manually resized columns, no longer auto "best fit")
$ws1.Columns.Item(3).ColumnWidth = 19.6
$ws1.Columns.Item(11).ColumnWidth = 36.3

# Selection moved to J2
$ws1.Range("J2").Select()

# ---------- Sheet "Notes" ----------
$ws2 = $wb.Worksheets.Item("Notes")

# Duplicate formatting of row 4 down into the newly used row 5 and row 6
$ws2.Range("F4:Q4").Copy()
$ws2.Range("F5:Q5").PasteSpecial(-4122)
$ws2.Range("F4:Q4").Copy()
$ws2.Range("F6:Q6").PasteSpecial(-4122)

$ws2.Range("E4").Copy()
$ws2.Range("E5").PasteSpecial(-4122)
$ws2.Range("E4").Copy()
$ws2.Range("E6").PasteSpecial(-4122)

# Content: split/expand the DialogueResponse note into three notes (rows 4,5,6)
$ws2.Range("E4").Value = "The response index starts from 0. If there are multiple instances of selection, the response index continues to increase."
$ws2.Range("E5").Value = "Ex: 1st instance with 5 responses, 2nd instance with 3 responses. To get the 2nd instance of index 2(last answer), the index would be 7."
$ws2.Range("E6").Value = "rewardKey < 0 add towards story progression. Empty rewardKey gives nothing."

# Selection moved to E5
$ws2.Range("E5").Select()
